# Handback report generation:
#  - "Handed back: in sync with en-US" -> "Handed back: not in sync with en-US"
#    (Status text shown on the Overview sheet and on each per-language sheet)
#  - New Correspond Handback DateTime stamps for the d19a2859... row in both
#    the zh-cn and de-de language sheets
#  - Widen the Status/zh-cn/de-de columns so the longer status text still fits
#    (mirrors the column autofit Excel performs after the text change)

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the new, longer text.
$overview.Range("E1").ColumnWidth = 32.62688700358076
$overview.Range("F1").ColumnWidth = 32.62688700358076

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# New handback report generated for the d19a2859... entry.
$zhcn.Range("K3").Value = "2016-09-06 05:23:27"

$zhcn.Range("C1").ColumnWidth = 32.62688700358076

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# New handback report generated for the d19a2859... entry.
$dede.Range("K3").Value = "2016-09-06 05:23:45"

$dede.Range("C1").ColumnWidth = 32.62688700358076
